# This script re-flows the two introductory paragraphs of the PowerShell
# tutorial into several runs each, one run per "visual line" of roughly
# constant width, with the separating space kept as its own run - matching
# a text-wrapped JSON/]export style.  It does this the same way a human
# editor driving real Word COM automation would: split each paragraph into
# a temporary paragraph-per-chunk sequence (clean, formatting-free splits),
# then delete the paragraph marks between the chunks so they collapse back
# into a single paragraph made up of multiple runs.

$d = $word.ActiveDocument

function Split-ParagraphIntoChunks($paraIndex, $offsets) {
    # $offsets are character positions (relative to the paragraph's Range.Start)
    # right before the separating space of each chunk boundary. For every
    # offset we insert two paragraph breaks: one right after the offset+1
    # (isolating the following text) and one right at the offset (isolating
    # the preceding text), so the single space in between becomes its own
    # paragraph. Processed back-to-front so earlier offsets stay valid.
    $p = $d.Paragraphs($paraIndex).Range
    $pstart = $p.Start
    for ($i = $offsets.Length - 1; $i -ge 0; $i--) {
        $off = $offsets[$i]

        $posAfterSpace = $pstart + $off + 1
        $rAfter = $d.Range($posAfterSpace, $posAfterSpace)
        $rAfter.InsertAfter([char]13)

        $posBeforeSpace = $pstart + $off
        $rBefore = $d.Range($posBeforeSpace, $posBeforeSpace)
        $rBefore.InsertAfter([char]13)
    }
}

function Merge-ParagraphsIntoRuns($paraIndex, $count) {
    # Collapses $count+1 consecutive paragraphs (starting at $paraIndex) back
    # into a single paragraph by deleting the paragraph mark at the end of
    # $paraIndex, $count times in a row. Each deletion keeps the text that
    # was in the following paragraph as a separate run (no run merging, no
    # leftover rPr), since it's a structural paragraph-mark deletion rather
    # than a text/content edit.
    for ($i = 0; $i -lt $count; $i++) {
        $r = $d.Paragraphs($paraIndex).Range
        $markStart = $r.End - 1
        $markEnd = $markStart + 1
        $markRange = $d.Range($markStart, $markEnd)
        $markRange.Delete()
    }
}

# --- Paragraph 2 ("FirstParagraph" style): the "Windows PowerShell is an
# evolution..." paragraph. Split it into 7 chunks (6 interior split points).
$offsets1 = @(70, 140, 213, 283, 356, 427)
Split-ParagraphIntoChunks 2 $offsets1

# --- Paragraph 3 ("BodyText" style): the "PowerShell is now an essential
# skill..." paragraph. Its original index shifted down because splitting
# paragraph 2 into 13 paragraphs (7 chunks + 6 spaces) added 12 new
# paragraphs ahead of it.
$para2Index = 3 + ($offsets1.Length * 2)
$offsets2 = @(70, 140, 213, 281, 350, 422)
Split-ParagraphIntoChunks $para2Index $offsets2

# Merge each run of temporary paragraphs back into a single paragraph.
$mergeCount1 = $offsets1.Length * 2
Merge-ParagraphsIntoRuns 2 $mergeCount1

$mergeCount2 = $offsets2.Length * 2
Merge-ParagraphsIntoRuns 3 $mergeCount2

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
